# Applies the profit-recalculation update from the scheduled-runner commit.
# Workbook "Sheets/Masamune_Profits.xlsx" holds one worksheet per crafting job
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); the runner refreshed the market-board
# price columns (H/I/J/K/L = prices, M/N = profit) for specific leve rows on each
# sheet. All cells here are plain cached numbers (no formulas in this workbook),
# so we just overwrite each changed cell with its new value via the COM object model.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 93
$ws.Cells.Item(93, 8).Value = 45998.4  # H93: 47598.4 -> 45998.4
$ws.Cells.Item(93, 10).Value = 45998.4  # J93: 47598.4 -> 45998.4
$ws.Cells.Item(93, 12).Value = 45998.4  # L93: 47598.4 -> 45998.4
$ws.Cells.Item(93, 14).Value = -50990.4  # N93: -52590.4 -> -50990.4
# row 95
$ws.Cells.Item(95, 8).Value = 37316.4  # H95: 37324.4 -> 37316.4
$ws.Cells.Item(95, 10).Value = 37316.4  # J95: 37324.4 -> 37316.4
$ws.Cells.Item(95, 12).Value = 37316.4  # L95: 37324.4 -> 37316.4
$ws.Cells.Item(95, 14).Value = -42808.4  # N95: -42816.4 -> -42808.4
# row 105
$ws.Cells.Item(105, 8).Value = 47996  # H105: 49335.5 -> 47996
$ws.Cells.Item(105, 10).Value = 47996  # J105: 49335.5 -> 47996
$ws.Cells.Item(105, 12).Value = 47996  # L105: 49335.5 -> 47996
$ws.Cells.Item(105, 14).Value = -54984  # N105: -56323.5 -> -54984
# row 123
$ws.Cells.Item(123, 8).Value = 37440.8  # H123: 37570 -> 37440.8
$ws.Cells.Item(123, 10).Value = 37440.8  # J123: 37570 -> 37440.8
$ws.Cells.Item(123, 12).Value = 37440.8  # L123: 37570 -> 37440.8
$ws.Cells.Item(123, 14).Value = -47240.8  # N123: -47370 -> -47240.8
# row 128
$ws.Cells.Item(128, 8).Value = 46766.668  # H128: 46776 -> 46766.668
$ws.Cells.Item(128, 10).Value = 46766.668  # J128: 46776 -> 46766.668
$ws.Cells.Item(128, 12).Value = 46766.668  # L128: 46776 -> 46766.668
$ws.Cells.Item(128, 14).Value = -56726.668  # N128: -56736 -> -56726.668

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 24
$ws.Cells.Item(24, 8).Value = 34915.668  # H24: 37803.332 -> 34915.668
$ws.Cells.Item(24, 10).Value = 34915.668  # J24: 37803.332 -> 34915.668
$ws.Cells.Item(24, 12).Value = 34915.668  # L24: 37803.332 -> 34915.668
$ws.Cells.Item(24, 14).Value = -35663.668  # N24: -38551.332 -> -35663.668
# row 63
$ws.Cells.Item(63, 8).Value = 3095.5625  # H63: 2419.5217 -> 3095.5625
$ws.Cells.Item(63, 9).Value = 2502.4167  # I63: 1953.0625 -> 2502.4167
$ws.Cells.Item(63, 10).Value = 4875  # J63: 3485.7144 -> 4875
$ws.Cells.Item(63, 11).Value = 2502.4167  # K63: 1953.0625 -> 2502.4167
$ws.Cells.Item(63, 12).Value = 4875  # L63: 3485.7144 -> 4875
$ws.Cells.Item(63, 13).Value = -1816.4167  # M63: -1267.0625 -> -1816.4167
$ws.Cells.Item(63, 14).Value = -6247  # N63: -4857.7144 -> -6247
# row 66
$ws.Cells.Item(66, 8).Value = 3095.5625  # H66: 2419.5217 -> 3095.5625
$ws.Cells.Item(66, 9).Value = 2502.4167  # I66: 1953.0625 -> 2502.4167
$ws.Cells.Item(66, 10).Value = 4875  # J66: 3485.7144 -> 4875
$ws.Cells.Item(66, 11).Value = 12512.0835  # K66: 9765.3125 -> 12512.0835
$ws.Cells.Item(66, 12).Value = 24375  # L66: 17428.572 -> 24375
$ws.Cells.Item(66, 13).Value = -9080.083500000001  # M66: -6333.3125 -> -9080.083500000001
$ws.Cells.Item(66, 14).Value = -31239  # N66: -24292.572 -> -31239
# row 95
$ws.Cells.Item(95, 8).Value = 40129.668  # H95: 40603.5 -> 40129.668
$ws.Cells.Item(95, 10).Value = 40129.668  # J95: 40603.5 -> 40129.668
$ws.Cells.Item(95, 12).Value = 40129.668  # L95: 40603.5 -> 40129.668
$ws.Cells.Item(95, 14).Value = -45621.668  # N95: -46095.5 -> -45621.668
# row 98
$ws.Cells.Item(98, 8).Value = 39785  # H98: 38115.668 -> 39785
$ws.Cells.Item(98, 10).Value = 39785  # J98: 38115.668 -> 39785
$ws.Cells.Item(98, 12).Value = 39785  # L98: 38115.668 -> 39785
$ws.Cells.Item(98, 14).Value = -45775  # N98: -44105.668 -> -45775
# row 100
$ws.Cells.Item(100, 8).Value = 34915.668  # H100: 37803.332 -> 34915.668
$ws.Cells.Item(100, 10).Value = 34915.668  # J100: 37803.332 -> 34915.668
$ws.Cells.Item(100, 12).Value = 34915.668  # L100: 37803.332 -> 34915.668
$ws.Cells.Item(100, 14).Value = -37079.668  # N100: -39967.332 -> -37079.668
# row 101
$ws.Cells.Item(101, 8).Value = 49594  # H101: 49582 -> 49594
$ws.Cells.Item(101, 10).Value = 49594  # J101: 49582 -> 49594
$ws.Cells.Item(101, 12).Value = 49594  # L101: 49582 -> 49594
$ws.Cells.Item(101, 14).Value = -56084  # N101: -56072 -> -56084
# row 103
$ws.Cells.Item(103, 8).Value = 42354  # H103: 42362 -> 42354
$ws.Cells.Item(103, 10).Value = 42354  # J103: 42362 -> 42354
$ws.Cells.Item(103, 12).Value = 42354  # L103: 42362 -> 42354
$ws.Cells.Item(103, 14).Value = -44698  # N103: -44706 -> -44698
# row 131
$ws.Cells.Item(131, 8).Value = 47803.5  # H131: 46897.75 -> 47803.5
$ws.Cells.Item(131, 10).Value = 47803.5  # J131: 46897.75 -> 47803.5
$ws.Cells.Item(131, 12).Value = 47803.5  # L131: 46897.75 -> 47803.5
$ws.Cells.Item(131, 14).Value = -57883.5  # N131: -56977.75 -> -57883.5
# row 137
$ws.Cells.Item(137, 8).Value = 43600  # H137: 44533.332 -> 43600
$ws.Cells.Item(137, 10).Value = 43600  # J137: 44533.332 -> 43600
$ws.Cells.Item(137, 12).Value = 43600  # L137: 44533.332 -> 43600
$ws.Cells.Item(137, 14).Value = -53800  # N137: -54733.332 -> -53800

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 92
$ws.Cells.Item(92, 8).Value = 39438.75  # H92: 39501.8 -> 39438.75
$ws.Cells.Item(92, 10).Value = 39438.75  # J92: 39501.8 -> 39438.75
$ws.Cells.Item(92, 12).Value = 39438.75  # L92: 39501.8 -> 39438.75
$ws.Cells.Item(92, 14).Value = -44430.75  # N92: -44493.8 -> -44430.75
# row 106
$ws.Cells.Item(106, 8).Value = 48663  # H106: 48331.5 -> 48663
$ws.Cells.Item(106, 10).Value = 48663  # J106: 48331.5 -> 48663
$ws.Cells.Item(106, 12).Value = 48663  # L106: 48331.5 -> 48663
$ws.Cells.Item(106, 14).Value = -51187  # N106: -50855.5 -> -51187
# row 130
$ws.Cells.Item(130, 8).Value = 46030.75  # H130: 47057.75 -> 46030.75
$ws.Cells.Item(130, 10).Value = 46030.75  # J130: 47057.75 -> 46030.75
$ws.Cells.Item(130, 12).Value = 46030.75  # L130: 47057.75 -> 46030.75
$ws.Cells.Item(130, 14).Value = -56070.75  # N130: -57097.75 -> -56070.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Cells.Item(16, 8).Value = 847.1539  # H16: 808.6923 -> 847.1539
$ws.Cells.Item(16, 10).Value = 934.3333  # J16: 767.6667 -> 934.3333
$ws.Cells.Item(16, 12).Value = 934.3333  # L16: 767.6667 -> 934.3333
$ws.Cells.Item(16, 14).Value = -1508.3333  # N16: -1341.6667 -> -1508.3333
# row 28
$ws.Cells.Item(28, 8).Value = 40321.5  # H28: 44000 -> 40321.5
$ws.Cells.Item(28, 10).Value = 40321.5  # J28: 44000 -> 40321.5
$ws.Cells.Item(28, 12).Value = 40321.5  # L28: 44000 -> 40321.5
$ws.Cells.Item(28, 14).Value = -40811.5  # N28: -44490 -> -40811.5
# row 43
$ws.Cells.Item(43, 8).Value = 48643  # H43: 48657 -> 48643
$ws.Cells.Item(43, 10).Value = 48643  # J43: 48657 -> 48643
$ws.Cells.Item(43, 12).Value = 48643  # L43: 48657 -> 48643
$ws.Cells.Item(43, 14).Value = -49011  # N43: -49025 -> -49011
# row 92
$ws.Cells.Item(92, 8).Value = 34312.125  # H92: 35199.75 -> 34312.125
$ws.Cells.Item(92, 10).Value = 34312.125  # J92: 35199.75 -> 34312.125
$ws.Cells.Item(92, 12).Value = 34312.125  # L92: 35199.75 -> 34312.125
$ws.Cells.Item(92, 14).Value = -39304.125  # N92: -40191.75 -> -39304.125
# row 96
$ws.Cells.Item(96, 8).Value = 37488.273  # H96: 43962.5 -> 37488.273
$ws.Cells.Item(96, 10).Value = 37488.273  # J96: 43962.5 -> 37488.273
$ws.Cells.Item(96, 12).Value = 37488.273  # L96: 43962.5 -> 37488.273
$ws.Cells.Item(96, 14).Value = -42980.273  # N96: -49454.5 -> -42980.273
# row 100
$ws.Cells.Item(100, 8).Value = 43413  # H100: 41618 -> 43413
$ws.Cells.Item(100, 10).Value = 43413  # J100: 41618 -> 43413
$ws.Cells.Item(100, 12).Value = 43413  # L100: 41618 -> 43413
$ws.Cells.Item(100, 14).Value = -45577  # N100: -43782 -> -45577
# row 101
$ws.Cells.Item(101, 8).Value = 48643  # H101: 48657 -> 48643
$ws.Cells.Item(101, 10).Value = 48643  # J101: 48657 -> 48643
$ws.Cells.Item(101, 12).Value = 48643  # L101: 48657 -> 48643
$ws.Cells.Item(101, 14).Value = -55133  # N101: -55147 -> -55133
# row 106
$ws.Cells.Item(106, 8).Value = 40424.168  # H106: 38978 -> 40424.168
$ws.Cells.Item(106, 10).Value = 48831.25  # J106: 49223.332 -> 48831.25
$ws.Cells.Item(106, 12).Value = 48831.25  # L106: 49223.332 -> 48831.25
$ws.Cells.Item(106, 14).Value = -51355.25  # N106: -51747.332 -> -51355.25
# row 113
$ws.Cells.Item(113, 8).Value = 847.1539  # H113: 808.6923 -> 847.1539
$ws.Cells.Item(113, 10).Value = 934.3333  # J113: 767.6667 -> 934.3333
$ws.Cells.Item(113, 12).Value = 934.3333  # L113: 767.6667 -> 934.3333
$ws.Cells.Item(113, 14).Value = -5274.3333  # N113: -5107.6667 -> -5274.3333
# row 124
$ws.Cells.Item(124, 8).Value = 35248  # H124: 36998 -> 35248
$ws.Cells.Item(124, 10).Value = 35248  # J124: 36998 -> 35248
$ws.Cells.Item(124, 12).Value = 35248  # L124: 36998 -> 35248
$ws.Cells.Item(124, 14).Value = -40158  # N124: -41908 -> -40158
# row 125
$ws.Cells.Item(125, 8).Value = 41162.5  # H125: 33495 -> 41162.5
$ws.Cells.Item(125, 10).Value = 41162.5  # J125: 33495 -> 41162.5
$ws.Cells.Item(125, 12).Value = 41162.5  # L125: 33495 -> 41162.5
$ws.Cells.Item(125, 14).Value = -46082.5  # N125: -38415 -> -46082.5
# row 131
$ws.Cells.Item(131, 8).Value = 41985  # H131: 41993 -> 41985
$ws.Cells.Item(131, 10).Value = 41985  # J131: 41993 -> 41985
$ws.Cells.Item(131, 12).Value = 41985  # L131: 41993 -> 41985
$ws.Cells.Item(131, 14).Value = -52065  # N131: -52073 -> -52065

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Cells.Item(80, 8).Value = 3855.9666  # H80: 3699.5 -> 3855.9666
$ws.Cells.Item(80, 9).Value = 4177.7856  # I80: 3824.625 -> 4177.7856
$ws.Cells.Item(80, 11).Value = 4177.7856  # K80: 3824.625 -> 4177.7856
$ws.Cells.Item(80, 13).Value = -3179.7856  # M80: -2826.625 -> -3179.7856
# row 83
$ws.Cells.Item(83, 8).Value = 3855.9666  # H83: 3699.5 -> 3855.9666
$ws.Cells.Item(83, 9).Value = 4177.7856  # I83: 3824.625 -> 4177.7856
$ws.Cells.Item(83, 11).Value = 20888.928  # K83: 19123.125 -> 20888.928
$ws.Cells.Item(83, 13).Value = -15896.928  # M83: -14131.125 -> -15896.928
# row 96
$ws.Cells.Item(96, 8).Value = 29965  # H96: 29969 -> 29965
$ws.Cells.Item(96, 10).Value = 29965  # J96: 29969 -> 29965
$ws.Cells.Item(96, 12).Value = 29965  # L96: 29969 -> 29965
$ws.Cells.Item(96, 14).Value = -35457  # N96: -35461 -> -35457
# row 105
$ws.Cells.Item(105, 8).Value = 0  # H105: 48671 -> 0
$ws.Cells.Item(105, 10).Value = 0  # J105: 48671 -> 0
$ws.Cells.Item(105, 12).Value = 0  # L105: 48671 -> 0
$ws.Cells.Item(105, 14).ClearContents()  # N105: -55659 -> (cell removed)
# row 127
$ws.Cells.Item(127, 8).Value = 26283.111  # H127: 26286 -> 26283.111
$ws.Cells.Item(127, 10).Value = 26283.111  # J127: 26286 -> 26283.111
$ws.Cells.Item(127, 12).Value = 26283.111  # L127: 26286 -> 26283.111
$ws.Cells.Item(127, 14).Value = -36203.111  # N127: -36206 -> -36203.111
# row 131
$ws.Cells.Item(131, 8).Value = 24997.75  # H131: 30772.666 -> 24997.75
$ws.Cells.Item(131, 10).Value = 24997.75  # J131: 30772.666 -> 24997.75
$ws.Cells.Item(131, 12).Value = 24997.75  # L131: 30772.666 -> 24997.75
$ws.Cells.Item(131, 14).Value = -35077.75  # N131: -40852.666 -> -35077.75
# row 137
$ws.Cells.Item(137, 8).Value = 51866.332  # H137: 46533 -> 51866.332
$ws.Cells.Item(137, 10).Value = 51866.332  # J137: 46533 -> 51866.332
$ws.Cells.Item(137, 12).Value = 51866.332  # L137: 46533 -> 51866.332
$ws.Cells.Item(137, 14).Value = -62066.332  # N137: -56733 -> -62066.332

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 106
$ws.Cells.Item(106, 8).Value = 35034  # H106: 42362 -> 35034
$ws.Cells.Item(106, 10).Value = 35034  # J106: 42362 -> 35034
$ws.Cells.Item(106, 12).Value = 35034  # L106: 42362 -> 35034
$ws.Cells.Item(106, 14).Value = -37558  # N106: -44886 -> -37558
# row 109
$ws.Cells.Item(109, 8).Value = 35217  # H109: 35239.668 -> 35217
$ws.Cells.Item(109, 10).Value = 35217  # J109: 35239.668 -> 35217
$ws.Cells.Item(109, 12).Value = 35217  # L109: 35239.668 -> 35217
$ws.Cells.Item(109, 14).Value = -37991  # N109: -38013.668 -> -37991
# row 117
$ws.Cells.Item(117, 8).Value = 43384  # H117: 43388 -> 43384
$ws.Cells.Item(117, 10).Value = 43384  # J117: 43388 -> 43384
$ws.Cells.Item(117, 12).Value = 43384  # L117: 43388 -> 43384
$ws.Cells.Item(117, 14).Value = -52562  # N117: -52566 -> -52562
# row 122
$ws.Cells.Item(122, 8).Value = 2575  # H122: 2516.6667 -> 2575
$ws.Cells.Item(122, 9).Value = 2525  # I122: 2420 -> 2525
$ws.Cells.Item(122, 10).Value = 2625  # J122: 3000 -> 2625
$ws.Cells.Item(122, 11).Value = 7575  # K122: 7260 -> 7575
$ws.Cells.Item(122, 12).Value = 7875  # L122: 9000 -> 7875
$ws.Cells.Item(122, 13).Value = -5125  # M122: -4810 -> -5125
$ws.Cells.Item(122, 14).Value = -12775  # N122: -13900 -> -12775
# row 123
$ws.Cells.Item(123, 8).Value = 39421  # H123: 39417 -> 39421
$ws.Cells.Item(123, 10).Value = 39421  # J123: 39417 -> 39421
$ws.Cells.Item(123, 12).Value = 39421  # L123: 39417 -> 39421
$ws.Cells.Item(123, 14).Value = -49221  # N123: -49217 -> -49221
# row 129
$ws.Cells.Item(129, 8).Value = 45418.332  # H129: 45118.57 -> 45418.332
$ws.Cells.Item(129, 10).Value = 45418.332  # J129: 45118.57 -> 45418.332
$ws.Cells.Item(129, 12).Value = 45418.332  # L129: 45118.57 -> 45418.332
$ws.Cells.Item(129, 14).Value = -55418.332  # N129: -55118.57 -> -55418.332
# row 130
$ws.Cells.Item(130, 8).Value = 48429  # H130: 48421 -> 48429
$ws.Cells.Item(130, 10).Value = 48429  # J130: 48421 -> 48429
$ws.Cells.Item(130, 12).Value = 48429  # L130: 48421 -> 48429
$ws.Cells.Item(130, 14).Value = -58469  # N130: -58461 -> -58469
# row 131
$ws.Cells.Item(131, 8).Value = 45097.332  # H131: 45318 -> 45097.332
$ws.Cells.Item(131, 10).Value = 45097.332  # J131: 45318 -> 45097.332
$ws.Cells.Item(131, 12).Value = 45097.332  # L131: 45318 -> 45097.332
$ws.Cells.Item(131, 14).Value = -55177.332  # N131: -55398 -> -55177.332
# row 139
$ws.Cells.Item(139, 8).Value = 50899.668  # H139: 49499.668 -> 50899.668
$ws.Cells.Item(139, 10).Value = 50899.668  # J139: 49499.668 -> 50899.668
$ws.Cells.Item(139, 12).Value = 50899.668  # L139: 49499.668 -> 50899.668
$ws.Cells.Item(139, 14).Value = -61179.668  # N139: -59779.668 -> -61179.668

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 94
$ws.Cells.Item(94, 8).Value = 40214.5  # H94: 40329 -> 40214.5
$ws.Cells.Item(94, 10).Value = 40214.5  # J94: 40329 -> 40214.5
$ws.Cells.Item(94, 12).Value = 40214.5  # L94: 40329 -> 40214.5
$ws.Cells.Item(94, 14).Value = -42016.5  # N94: -42131 -> -42016.5
# row 109
$ws.Cells.Item(109, 8).Value = 39369  # H109: 39373 -> 39369
$ws.Cells.Item(109, 10).Value = 39369  # J109: 39373 -> 39369
$ws.Cells.Item(109, 12).Value = 39369  # L109: 39373 -> 39369
$ws.Cells.Item(109, 14).Value = -42143  # N109: -42147 -> -42143
# row 118
$ws.Cells.Item(118, 8).Value = 43380  # H118: 43372 -> 43380
$ws.Cells.Item(118, 10).Value = 43380  # J118: 43372 -> 43380
$ws.Cells.Item(118, 12).Value = 43380  # L118: 43372 -> 43380
$ws.Cells.Item(118, 14).Value = -46694  # N118: -46686 -> -46694
# row 127
$ws.Cells.Item(127, 8).Value = 42421  # H127: 42429 -> 42421
$ws.Cells.Item(127, 10).Value = 42421  # J127: 42429 -> 42421
$ws.Cells.Item(127, 12).Value = 42421  # L127: 42429 -> 42421
$ws.Cells.Item(127, 14).Value = -52341  # N127: -52349 -> -52341

Write-Host "Updated 208 cell(s) and cleared 1 cell(s) across 7 sheets."
